# Rename the header columns to reflect the respective input file's format
# version (FV2310 / FV2404) instead of the generic "_old" / "_new" suffixes,
# then turn the sheet's used range into a real Excel Table and freeze the
# header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the shared-string header labels: "_old" -> "_FV2310", "_new" -> "_FV2404"
$headers = @(
    "Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID",
    "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung"
)

for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = "$($headers[$col - 1])_FV2310"
}

for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = "$($headers[$col - 12])_FV2404"
}

# 2) Convert the data range into an Excel Table ("Table1")
$tableRange = $ws.Range("A1:U72")
$table = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$table.Name = "Table1"

# 3) Freeze the header row (selection must sit on row 2 for the freeze
#    boundary to land between row 1 and row 2), then restore the selection
#    back to A1 to match the default/no-explicit-selection appearance.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
